$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 353, shifting existing rows 353..449 down to 354..450
# (dimension grows from A1:R449 to A1:R450).
$ws.Rows.Item(353).Insert()

# Populate the newly inserted row 353 with the new data record.
# Columns A,B,C,E,F,G,H,I,R carry the same constant category metadata as
# every other data row on this sheet; D,J,K,L,M,P are the new values; N,O,Q
# match the values that were already present for this record.
$ws.Range("A353").Value = 10
$ws.Range("B353").Value = "Vega Modelo de Temuco"
$ws.Range("C353").Value = "La Araucanía"
$ws.Range("D353").Value = 44798
$ws.Range("E353").Value = 9
$ws.Range("F353").Value = 100114014
$ws.Range("G353").Value = "Betarraga"
$ws.Range("H353").Value = "Sin especificar"
$ws.Range("I353").Value = "Primera"
$ws.Range("J353").Value = 130
$ws.Range("K353").Value = 10000
$ws.Range("L353").Value = 11000
$ws.Range("M353").Value = 10385
$ws.Range("N353").Value = "$/docena de paquetes"
$ws.Range("O353").Value = "Provincia de Cautín"
$ws.Range("P353").Value = 865
$ws.Range("Q353").Value = 12
$ws.Range("R353").Value = "Hortaliza"
